# Applies the cryptos.xlsx price/volume/ranking update described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.861.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.10%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.085.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.38%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'233.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.05%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.26%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'59.26"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.85%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.395"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.53%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0788"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.64%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.37%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.391.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.19%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'14.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'21.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.18%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +1.38%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'5.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.91%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.155.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.95%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'37.785.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.14%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.25%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'71.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.97%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0848"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.02%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'227.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.28%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E24").Value = "'  -0.59%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.86%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'171.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.40%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.50%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.58%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.41%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'19.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.08%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.54%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.22%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0633"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.65%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.21%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.95%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.04%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.71%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.12%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.50%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.32%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'99.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.02%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'17.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +8.38%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0218"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.95%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.55%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.449.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.26%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.82%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +2.89%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.49%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'7.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.40%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.76%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.276.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.39%  "
$ws.Range("E51").Style = "Normal"
